# "Add dataset 9 to openjml output"
#
# dataset_id 9 (rows 13-16) previously had too few judged snippets with
# warnings for openjml, so num_snippets_warnings/num_warnings were 0 and no
# correlation stats were computed on the "openjml" sheet. This commit adds
# the openjml warning/correlation numbers for dataset 9, which also updates
# the aggregated numbers (and derived correlations) on the "all_tools" sheet.

$wb = $excel.ActiveWorkbook

# ---- "all_tools" sheet: refresh dataset 9 (rows 13-16) now that openjml
#      contributes warning counts/correlations for it too. ----
$wsAll = $wb.Worksheets.Item("all_tools")

$wsAll.Range("F13").Value = 30
$wsAll.Range("G13").Value = 108
$wsAll.Range("I13").Value = -0.3546780165447971
$wsAll.Range("J13").Value = 0.01262159666689494
$wsAll.Range("K13").Value = -0.4463770963226328
$wsAll.Range("L13").Value = 0.01341282859891604

$wsAll.Range("F14").Value = 30
$wsAll.Range("G14").Value = 108
$wsAll.Range("I14").Value = 0.2846501884412231
$wsAll.Range("J14").Value = 0.04660031288457147
$wsAll.Range("K14").Value = 0.3544834006793617
$wsAll.Range("L14").Value = 0.05460272238348809

$wsAll.Range("F15").Value = 30
$wsAll.Range("G15").Value = 108
$wsAll.Range("I15").Value = 0.2886751345948129
$wsAll.Range("J15").Value = 0.04276947240054693
$wsAll.Range("K15").Value = 0.3717143584291471
$wsAll.Range("L15").Value = 0.04311672358740608

$wsAll.Range("F16").Value = 30
$wsAll.Range("G16").Value = 108
$wsAll.Range("I16").Value = -0.3462790510727774
$wsAll.Range("J16").Value = 0.01478142552429735
$wsAll.Range("K16").Value = -0.4190867990009841
$wsAll.Range("L16").Value = 0.02115728721848806

# ---- "openjml" sheet: add dataset 9's warning counts/correlation stats. ----
$wsOpenjml = $wb.Worksheets.Item("openjml")

$wsOpenjml.Range("F13").Value = 30
$wsOpenjml.Range("G13").Value = 78
$wsOpenjml.Range("I13").Value = -0.1306708482007147
$wsOpenjml.Range("J13").Value = 0.3581258439224645
$wsOpenjml.Range("K13").Value = -0.2011840434130176
$wsOpenjml.Range("L13").Value = 0.2863985630278126

$wsOpenjml.Range("F14").Value = 30
$wsOpenjml.Range("G14").Value = 78
$wsOpenjml.Range("I14").Value = 0.04028068704356932
$wsOpenjml.Range("J14").Value = 0.7782565834710006
$wsOpenjml.Range("K14").Value = 0.06082851254261235
$wsOpenjml.Range("L14").Value = 0.7494911940114868

$wsOpenjml.Range("F15").Value = 30
$wsOpenjml.Range("G15").Value = 78
$wsOpenjml.Range("I15").Value = 0.08553337321327789
$wsOpenjml.Range("J15").Value = 0.5483178177462631
$wsOpenjml.Range("K15").Value = 0.1166847704091495
$wsOpenjml.Range("L15").Value = 0.5391786637186142

$wsOpenjml.Range("F16").Value = 30
$wsOpenjml.Range("G16").Value = 78
$wsOpenjml.Range("I16").Value = -0.2077674306436665
$wsOpenjml.Range("J16").Value = 0.1435729559776537
$wsOpenjml.Range("K16").Value = -0.2975516272906987
$wsOpenjml.Range("L16").Value = 0.1102840897457505
